# feat: add 2022-Q1 data
#
# - Adds a new "2022-Q1" worksheet (fund holdings), positioned between the
#   existing "2021-Q4" sheet and the "总计" (totals) sheet.
# - Adds a new leading row to the "总计" sheet summarising the 2022-Q1 data,
#   pushing the existing 2021-Q4 total row down.

$wb = $excel.ActiveWorkbook

# Helper: force a numeric-looking value to be stored as TEXT (matches the
# source data, where columns like "基金规模"/"仓位占比" are text, not numbers),
# then drop back to the default "Normal" style so we don't leave a stray
# number-format style applied to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$q4 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted right after "2021-Q4".
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Match the outline / page-setup style of the sibling sheets.
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72

# Start from the "2021-Q4" layout/formatting, then overwrite with the new data.
$q4.Range("A1:H3").Copy($q1.Range("A1"))

Set-TextValue $q1.Range("B2") "004854"
$q1.Range("C2").Value = "广发中证全指汽车指数A"
Set-TextValue $q1.Range("D2") "22.01"
Set-TextValue $q1.Range("E2") "94.43"
Set-TextValue $q1.Range("F2") "3.41"
Set-TextValue $q1.Range("G2") "0.7505"
$q1.Range("H2").Value = 8

Set-TextValue $q1.Range("B3") "004855"
$q1.Range("C3").Value = "广发中证全指汽车指数C"
Set-TextValue $q1.Range("D3") "6.11"
Set-TextValue $q1.Range("E3") "94.43"
Set-TextValue $q1.Range("F3") "3.41"
Set-TextValue $q1.Range("G3") "0.2084"
$q1.Range("H3").Value = 8

# ---------------------------------------------------------------------
# 2. "总计" sheet: add a 2022-Q1 total row above the existing 2021-Q4 row.
#    Re-resolve the sheet by name now that a sheet has been inserted, so we
#    pick up the worksheet itself rather than a stale positional reference.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing (2021-Q4) row down to row 3, carrying its formatting,
# then overwrite row 2 in place with the new 2022-Q1 totals.
$total.Range("A2:D2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1

$total.Range("A2").Value = 0
Set-TextValue $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.96
